$d = $word.ActiveDocument

# Locate the "Dependencies" heading paragraph and the paragraph that
# documents the Visual C++ 2012 Redistributable dependency, then the
# heading paragraph that follows them ("How to develop for Redis").
$depHeadingIndex = -1
$depBodyIndex = -1
$nextHeadingIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "^Dependencies\s*$") {
        $depHeadingIndex = $i
    } elseif ($depHeadingIndex -ne -1 -and $depBodyIndex -eq -1 -and $t -match "dependent on the") {
        $depBodyIndex = $i
    } elseif ($depBodyIndex -ne -1 -and $nextHeadingIndex -eq -1) {
        $nextHeadingIndex = $i
        break
    }
}

# Remove the "Dependencies" heading paragraph and the Visual C++
# redistributable paragraph that follows it (the release no longer
# documents the DLL CRT redistributable dependency).
$depHeading = $d.Paragraphs.Item($depHeadingIndex)
$depBody = $d.Paragraphs.Item($depBodyIndex)
$rng = $d.Range($depHeading.Range.Start, $depBody.Range.End)
$rng.Delete()

# Re-anchor the _GoBack bookmark (previously trailing the "Running as a
# Service" paragraph) onto the start of the next heading paragraph,
# which is now "How to develop for Redis".
$nextHeading = $d.Paragraphs.Item($depHeadingIndex)
$bmRange = $d.Range($nextHeading.Range.Start, $nextHeading.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
